$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as plain text so that values such as
# "1.000", "4.970", "0.00001028", "8.500", "0.9380", "0.1140" are not
# reinterpreted by Excel as numbers (which would silently drop trailing
# zeros / switch to scientific notation and lose the original formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.347.03"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.841.52"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "239.93"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "0.6294"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.07442"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "24.97"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").Value = "0.2893"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.841.05"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "4.970"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "0.6755"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "0.00001028"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "81.79"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "6.237"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "29.354.10"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "228.62"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "12.32"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "7.370"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "158.07"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "8.500"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").Value = "17.44"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "0.06913"
$ws.Range("E28").Value = "  +6.90%  "
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").Value = "1.486"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "4.058"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "1.822"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "1.138"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "0.6983"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "2.585"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D38").Value = "2.821"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").Value = "1.237.72"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "6.807"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("D41").Value = "0.9380"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("D42").Value = "0.9991"
$ws.Range("D43").Value = "1.989.39"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "65.28"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("E46").Value = "  +3.79%  "
$ws.Range("D47").Value = "7.022"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "1.707"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "8.977"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "0.1140"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "0.3907"
$ws.Range("E51").Value = "  -1.05%  "
